$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 91 and 92 (fixtures 6924568 / 6924569) had their record order swapped:
# what used to be row 91 (Atletico Morelia vs Atlante) is now row 92, and what used to be
# row 92 (Venados FC vs Dorados) is now row 91. Swap every data column (B:AB).

# New row 91 <= old row 92
$ws.Range("B91").Value = 6924569
$ws.Range("C91").Value = "Mexico Liga de Expansion"
$ws.Range("D91").Value = 45214.92013888889
$ws.Range("E91").Value = "Venados FC"
$ws.Range("F91").Value = "Dorados"
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 1
$ws.Range("I91").Value = "H"
$ws.Range("J91").Value = 1.615
$ws.Range("K91").Value = 4
$ws.Range("L91").Value = 4.5
$ws.Range("M91").Value = 1.5
$ws.Range("N91").Value = 4.75
$ws.Range("O91").Value = 5.75
$ws.Range("P91").Value = -1.25
$ws.Range("Q91").Value = 1.925
$ws.Range("R91").Value = 1.875
$ws.Range("S91").Value = 3
$ws.Range("T91").Value = 1.75
$ws.Range("U91").Value = 1.95
$ws.Range("V91").Value = 0.5
$ws.Range("W91").Value = -1
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = 0.925
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.75
$ws.Range("AB91").Value = -1

# New row 92 <= old row 91
$ws.Range("B92").Value = 6924568
$ws.Range("C92").Value = "Mexico Liga de Expansion"
$ws.Range("D92").Value = 45214.92013888889
$ws.Range("E92").Value = "Atletico Morelia"
$ws.Range("F92").Value = "Atlante"
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 1
$ws.Range("I92").Value = "A"
$ws.Range("J92").Value = 2.4
$ws.Range("K92").Value = 3
$ws.Range("L92").Value = 2.875
$ws.Range("M92").Value = 2.7
$ws.Range("N92").Value = 3.1
$ws.Range("O92").Value = 2.8
$ws.Range("P92").Value = 0
$ws.Range("Q92").Value = 1.85
$ws.Range("R92").Value = 1.95
$ws.Range("S92").Value = 2.25
$ws.Range("T92").Value = 1.975
$ws.Range("U92").Value = 1.725
$ws.Range("V92").Value = -1
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = 1.8
$ws.Range("Y92").Value = -1
$ws.Range("Z92").Value = 0.95
$ws.Range("AA92").Value = -1
$ws.Range("AB92").Value = 0.7250000000000001

# --- Rows 186 and 187 (fixtures 7648957 / 7648958) also had their match-id/odds swapped
# between the two rows, while the two fixtures' team names (Monterrey U23 vs Mazatlan FC U23,
# and Unam Pumas U23 vs Tijuana U23) stay on the same rows. Swap B and G:AB, leave E/F alone.

# New row 186 <= old row 187 (B, G:AB)
$ws.Range("B186").Value = 7648957
$ws.Range("C186").Value = "Mexico Liga de Expansion"
$ws.Range("D186").Value = 45361.5
$ws.Range("G186").Value = 2
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = "H"
$ws.Range("J186").Value = 1.666
$ws.Range("K186").Value = 3.5
$ws.Range("L186").Value = 4.2
$ws.Range("M186").Value = 1.533
$ws.Range("N186").Value = 4.333
$ws.Range("O186").Value = 6
$ws.Range("P186").Value = -1.25
$ws.Range("Q186").Value = 2.025
$ws.Range("R186").Value = 1.775
$ws.Range("S186").Value = 2.75
$ws.Range("T186").Value = 1.775
$ws.Range("U186").Value = 2.025
$ws.Range("V186").Value = 0.5329999999999999
$ws.Range("W186").Value = -1
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = 1.025
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 1.025

# New row 187 <= old row 186 (B, G:AB)
$ws.Range("B187").Value = 7648958
$ws.Range("C187").Value = "Mexico Liga de Expansion"
$ws.Range("D187").Value = 45361.5
$ws.Range("G187").Value = 4
$ws.Range("H187").Value = 3
$ws.Range("I187").Value = "H"
$ws.Range("J187").Value = 2.375
$ws.Range("K187").Value = 3.1
$ws.Range("L187").Value = 2.75
$ws.Range("M187").Value = 2.375
$ws.Range("N187").Value = 3.4
$ws.Range("O187").Value = 3
$ws.Range("P187").Value = -0.25
$ws.Range("Q187").Value = 2
$ws.Range("R187").Value = 1.8
$ws.Range("S187").Value = 2.75
$ws.Range("T187").Value = 1.95
$ws.Range("U187").Value = 1.85
$ws.Range("V187").Value = 1.375
$ws.Range("W187").Value = -1
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = 1
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 0.95
$ws.Range("AB187").Value = -1
